$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Depth First Search")

# Seed rows 9-12 by duplicating the formatting of the last existing row (row 8).
$ws.Range("A8:K8").Copy()
$ws.Range("A9:K12").PasteSpecial(-4122)

# Row heights: rows 9-10 match the plain rows (20), rows 11-12 need extra
# height (30) because column E will hold a wrapped multi-line note.
$ws.Range("A9:K10").RowHeight = 20
$ws.Range("A11:K12").RowHeight = 30

# --- Row 9: Find Largest Value in Each Tree Row -----------------------
$ws.Range("A9").Value = 515
$ws.Range("B9").Value = "Find Largest Value in Each Tree Row"
$ws.Range("C9").Value = "https://leetcode.com/problems/find-largest-value-in-each-tree-row/description/"
$ws.Range("D9").Value = "Medium"
$ws.Range("F9").Value = "Binary Tree"
$ws.Range("G9").Value = "O(n)"
$ws.Range("H9").Value = "O(logn)"
$ws.Range("I9").Value = 45488

# --- Row 10: Add One Row to Tree ---------------------------------------
$ws.Range("A10").Value = 623
$ws.Range("B10").Value = "Add One Row to Tree"
$ws.Range("C10").Value = "https://leetcode.com/problems/add-one-row-to-tree/description/"
$ws.Range("D10").Value = "Medium"
$ws.Range("F10").Value = "Binary Tree"
$ws.Range("G10").Value = "O(n)"
$ws.Range("H10").Value = "O(1)"
$ws.Range("I10").Value = 45488

# --- Row 11: Find Duplicate Subtrees -----------------------------------
$ws.Range("A11").Value = 652
$ws.Range("B11").Value = "Find Duplicate Subtrees"
$ws.Range("C11").Value = "https://leetcode.com/problems/find-duplicate-subtrees/description/"
$ws.Range("D11").Value = "Medium"
$ws.Range("E11").Value = "When marking null child nodes,`nuse a slash to help denote left or right."
$ws.Range("E11").WrapText = $true
$ws.Range("F11").Value = "Binary Tree"
$ws.Range("G11").Value = "O(n)"
$ws.Range("H11").Value = "O(n)"
$ws.Range("I11").Value = 45488

# --- Row 12: Longest Univalue Path -------------------------------------
$ws.Range("A12").Value = 687
$ws.Range("B12").Value = "Longest Univalue Path"
$ws.Range("C12").Value = "https://leetcode.com/problems/longest-univalue-path/description/"
$ws.Range("D12").Value = "Medium"
$ws.Range("E12").Value = "Beware of ""semi-cycles"" in longest `nunivalue paths."
$ws.Range("E12").WrapText = $true
$ws.Range("F12").Value = "Binary Tree"
$ws.Range("G12").Value = "O(n)"
$ws.Range("H12").Value = "O(1)"
$ws.Range("I12").Value = 45488

# Real hyperlinks (relationship + <hyperlink> element) for the new URLs.
$ws.Hyperlinks.Add($ws.Range("C9"), "https://leetcode.com/problems/find-largest-value-in-each-tree-row/description/")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://leetcode.com/problems/add-one-row-to-tree/description/")
$ws.Hyperlinks.Add($ws.Range("C11"), "https://leetcode.com/problems/find-duplicate-subtrees/description/")
$ws.Hyperlinks.Add($ws.Range("C12"), "https://leetcode.com/problems/longest-univalue-path/description/")

# Hyperlinks.Add re-styles the anchor cell with a generic hyperlink style;
# restore the sheet's own hyperlink-cell format (same as C2:C8) on column C.
$ws.Range("C8").Copy()
$ws.Range("C9:C12").PasteSpecial(-4122)

# Extend the Easy/Medium/Hard conditional formatting down to the new rows.
$ws.Range("D2:D8").FormatConditions.Item(1).ModifyAppliesToRange($ws.Range("D2:D12"))

$excel.CutCopyMode = $false
